$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 33) with the new user's details
$ws.Range("A33").Value = 110032
$ws.Range("B33").Value = 9317596770
$ws.Range("C33").Value = "Ewan Marsh"
$ws.Range("D33").Value = "ewan.marsh@xyz.com"
$ws.Range("E33").Value = 818876433
$ws.Range("F33").Value = "ACT"
$ws.Range("G33").Value = "eng"
$ws.Range("H33").Value = "PWD"
$ws.Range("I33").Value = $true
$ws.Range("J33").Value = "superadmin"
$ws.Range("K33").Value = "now()"
$ws.Range("L33").Value = "now()"

$wb.Save()
